$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value to 2
$ws.Range("B2").Value = 2

# Delete row 3 entirely
$ws.Rows(3).Delete()
